$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "RepaymentStrategy" label (row 18) is now lower-cased to "repaymentstrategy"
$ws.Range("A18").Value = "repaymentstrategy"

# Move the sheet's active selection to reflect where the author was working
$ws.Range("C18").Select()
